$wb = $excel.ActiveWorkbook

# --- 1) SCALE_CONSISTENCY sheet: insert two columns before column B ---
# This shifts existing data from B:F to D:H, matching the new dimension A1:H12.
$wsConsistency = $wb.Worksheets.Item("SCALE_CONSISTENCY")
$wsConsistency.Columns("B:C").Insert()

# --- 2) CONFIDENCE_INTERVALS sheet: replace computed numeric results with
#        #NAME? errors in the "Confidence" / "Confidence interval" columns
#        (both the per-item block E:G and the per-scale block M:O). ---
$wsConfidence = $wb.Worksheets.Item("CONFIDENCE_INTERVALS")

$rows = 5..12
foreach ($r in $rows) {
    $wsConfidence.Cells.Item($r, 5).Value = "#NAME?"   # column E
    $wsConfidence.Cells.Item($r, 6).Value = "#NAME?"   # column F
    $wsConfidence.Cells.Item($r, 7).Value = "#NAME?"   # column G
}

$rowsScale = 5..7
foreach ($r in $rowsScale) {
    $wsConfidence.Cells.Item($r, 13).Value = "#NAME?"  # column M
    $wsConfidence.Cells.Item($r, 14).Value = "#NAME?"  # column N
    $wsConfidence.Cells.Item($r, 15).Value = "#NAME?"  # column O
}
